# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The G column ("K") values were recalculated from source data (K count
# instead of the previous Strike# metric) and rewritten into the sheet.
# All other columns/rows are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-27 (corresponds to G2:G27 in the sheet)
$newK = @(4, 5, 7, 3, 3, 1, 4, 3, 6, 6, 6, 9, 4, 4, 4, 5, 6, 2, 5, 4, 6, 3, 6, 2, 2, 1)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
